# Applies the "a paper on multiphase SPH" commit:
#  1) Merge the split bold runs of the "Oil-Spill Simulation Using Bi-Layer
#     Shallow Water equations:" heading into a single run.
#  2) Merge the split bold runs of the "Real-time Collision Detection and
#     Distance Computation on Point Cloud Sensor Data" heading into a single
#     run.
#  3) Replace the empty paragraph that used to hold only the `_GoBack`
#     bookmark (after the "Particle-Based Fluid Simulation..." entry) with a
#     new list entry ("Multiphase SPH Simulation for Interactive Fluids and
#     Solids:") plus its summary paragraph, re-homing the bookmark inside the
#     new summary text.

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-RangeXml($range, [string]$innerXml) {
    $xml = $pkgOpen + $innerXml + $pkgClose
    $range.InsertXML($xml)
}

function Find-ParagraphIndex($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    throw "Paragraph containing '$needle' not found"
}

# Resolve every paragraph index up front: none of the edits below add or
# remove paragraphs (step 1/2 merge runs in place, step 3 swaps the contents
# of a single existing paragraph for a bigger block), so indices stay valid
# for the whole script.
$oilIdx = Find-ParagraphIndex $d "Oil-Spill Simulation Using"
$rtIdx = Find-ParagraphIndex $d "Real-time Collision Detection"
$particleIdx = Find-ParagraphIndex $d "Particle-Based Fluid Simulation for Interactive Applications"
$summaryIdx = $particleIdx + 1
$bookmarkIdx = $particleIdx + 2

# --- 1) "Oil-Spill Simulation Using Bi-Layer Shallow Water equations:" ---
$oilXml = '<w:p><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-GB"/></w:rPr><w:t>Oil-Spill Simulation Using Bi-Layer Shallow Water equations:</w:t></w:r></w:p>'
Set-RangeXml $d.Paragraphs.Item($oilIdx).Range $oilXml

# --- 2) "Real-time Collision Detection and Distance Computation on Point Cloud Sensor Data" ---
$rtXml = '<w:p><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-GB"/></w:rPr><w:t>Real-time Collision Detection and Distance Computation on Point Cloud Sensor Data</w:t></w:r></w:p>'
Set-RangeXml $d.Paragraphs.Item($rtIdx).Range $rtXml

# --- 3) Replace the bookmark-only paragraph with the new Multiphase SPH entry ---
$newBlock = @'
<w:p><w:pPr><w:ind w:left="720"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Lijstalinea"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-GB"/></w:rPr><w:t>Multiphase SPH Simulation for Interactive Fluids and Solids</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="en-GB"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Lijstalinea"/><w:rPr><w:b/><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Lijstalinea"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">This paper introduces solid phases, </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>deformable bodies and granular materials</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> to the multiphase SPH. With their implementation on the GPU they </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>can</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> set their timestep to something between 10</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/><w:lang w:val="en-GB"/></w:rPr><w:t>-3</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> and 10</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/><w:lang w:val="en-GB"/></w:rPr><w:t>-4</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>. They also reference some other papers which work with the SPH framework. We don’t need the whole multiphase part, but it should be possible to use some ideas of the SPH GPU implementation.</w:t></w:r></w:p>
'@

Set-RangeXml $d.Paragraphs.Item($bookmarkIdx).Range $newBlock
